$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.981.34'
$ws.Range('E2').Value = '  -2.23%  '
$ws.Range('D3').Value = '2.256.90'
$ws.Range('E3').Value = '  -3.55%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '297.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.497'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.64%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -3.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.87'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.81%  '
$ws.Range('E11').Value = '  -1.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -9.21%  '
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').Value = '2.607.60'
$ws.Range('E15').Value = '  -3.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.83%  '
$ws.Range('D17').Value = '2.256.62'
$ws.Range('E17').Value = '  -5.96%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.772'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.12%  '
$ws.Range('D19').Value = '42.001.47'
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('D20').Value = '0.0₃0887'
$ws.Range('E20').Value = '  -2.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.33'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.48'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '232.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('E25').Value = '  -4.78%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.44'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.65%  '
$ws.Range('E29').Value = '  -7.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.36'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.34'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0691'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.01%  '
$ws.Range('E38').Value = '  -6.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.94%  '
$ws.Range('E40').Value = '  -4.52%  '
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.77%  '
$ws.Range('E43').Value = '  +1.73%  '
$ws.Range('D44').Value = '1.942.72'
$ws.Range('E44').Value = '  -3.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0277'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.02%  '
$ws.Range('E47').Value = '  -7.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.81'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.06%  '
$ws.Range('D50').Value = '2.482.07'
$ws.Range('E50').Value = '  -3.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.54%  '
